$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H45").Value = 0
$ws.Range("I45").Value = 0
$ws.Range("K45").Value = 0
$ws.Range("M45").ClearContents()
$ws.Range("H76").Value = 6497.1113
$ws.Range("I76").Value = 7747.5
$ws.Range("K76").Value = 7747.5
$ws.Range("M76").Value = -7432.5
$ws.Range("H79").Value = 6497.1113
$ws.Range("I79").Value = 7747.5
$ws.Range("K79").Value = 7747.5
$ws.Range("M79").Value = -6655.5
$ws.Range("H112").Value = 53048.8
$ws.Range("J112").Value = 55821.105
$ws.Range("L112").Value = 167463.315
$ws.Range("N112").Value = -169679.315
$ws.Range("H116").Value = 9331.308000000001
$ws.Range("I116").Value = 10659.333
$ws.Range("J116").Value = 8193
$ws.Range("K116").Value = 10659.333
$ws.Range("L116").Value = 8193
$ws.Range("M116").Value = -7217.333000000001
$ws.Range("N116").Value = -15077
$ws.Range("H132").Value = 1078.4474
$ws.Range("I132").Value = 1064.909
$ws.Range("K132").Value = 3194.727
$ws.Range("M132").Value = -664.7270000000003
$ws.Range("H138").Value = 3905.9111
$ws.Range("I138").Value = 2101.077
$ws.Range("J138").Value = 4639.125
$ws.Range("K138").Value = 6303.231000000001
$ws.Range("L138").Value = 13917.375
$ws.Range("M138").Value = -1163.231000000001
$ws.Range("N138").Value = -24197.375

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H69").Value = 500459
$ws.Range("J69").Value = 500459
$ws.Range("L69").Value = 500459
$ws.Range("N69").Value = -501957
$ws.Range("H72").Value = 500459
$ws.Range("J72").Value = 500459
$ws.Range("L72").Value = 1501377
$ws.Range("N72").Value = -1508865
$ws.Range("H132").Value = 2627.7778
$ws.Range("I132").Value = 2142.484
$ws.Range("K132").Value = 6427.451999999999
$ws.Range("M132").Value = -3897.451999999999

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H81").Value = 18250
$ws.Range("J81").Value = 18250
$ws.Range("L81").Value = 18250
$ws.Range("N81").Value = -20372
$ws.Range("H84").Value = 18250
$ws.Range("J84").Value = 18250
$ws.Range("L84").Value = 54750
$ws.Range("N84").Value = -65358
$ws.Range("H107").Value = 2376.081
$ws.Range("I107").Value = 2709.5715
$ws.Range("J107").Value = 1338.5555
$ws.Range("K107").Value = 2709.5715
$ws.Range("L107").Value = 1338.5555
$ws.Range("M107").Value = -789.5715
$ws.Range("N107").Value = -5178.5555
$ws.Range("H128").Value = 10020.5
$ws.Range("I128").Value = 10020.5
$ws.Range("K128").Value = 30061.5
$ws.Range("M128").Value = -27571.5

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 3200.3333
$ws.Range("I7").Value = 1023
$ws.Range("K7").Value = 1023
$ws.Range("M7").Value = -910
$ws.Range("H31").Value = 119496.336
$ws.Range("I31").Value = 169416.33
$ws.Range("J31").Value = 19656.334
$ws.Range("K31").Value = 169416.33
$ws.Range("L31").Value = 19656.334
$ws.Range("M31").Value = -169121.33
$ws.Range("N31").Value = -20246.334
$ws.Range("H34").Value = 119496.336
$ws.Range("I34").Value = 169416.33
$ws.Range("J34").Value = 19656.334
$ws.Range("K34").Value = 169416.33
$ws.Range("L34").Value = 19656.334
$ws.Range("M34").Value = -169214.33
$ws.Range("N34").Value = -20060.334
$ws.Range("H58").Value = 2652.2
$ws.Range("I58").Value = 2504.6924
$ws.Range("K58").Value = 2504.6924
$ws.Range("M58").Value = -2301.6924
$ws.Range("H132").Value = 2873
$ws.Range("I132").Value = 2825.6943
$ws.Range("J132").Value = 3440.6667
$ws.Range("K132").Value = 8477.082900000001
$ws.Range("L132").Value = 10322.0001
$ws.Range("M132").Value = -5947.082900000001
$ws.Range("N132").Value = -15382.0001
$ws.Range("H134").Value = 9808.029
$ws.Range("I134").Value = 6759.56
$ws.Range("K134").Value = 20278.68
$ws.Range("M134").Value = -17743.68
$ws.Range("H136").Value = 2652.2
$ws.Range("I136").Value = 2504.6924
$ws.Range("K136").Value = 7514.0772
$ws.Range("M136").Value = -4964.0772
$ws.Range("H137").Value = 92936.25
$ws.Range("J137").Value = 92936.25
$ws.Range("L137").Value = 92936.25
$ws.Range("N137").Value = -103136.25

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 3148881.5
$ws.Range("I4").Value = 3412405
$ws.Range("K4").Value = 10237215
$ws.Range("M4").Value = -10237103
$ws.Range("H38").Value = 105.42857
$ws.Range("I38").Value = 62.5
$ws.Range("K38").Value = 187.5
$ws.Range("M38").Value = 159.5
$ws.Range("H107").Value = 720.6667
$ws.Range("J107").Value = 946
$ws.Range("L107").Value = 2838
$ws.Range("N107").Value = -6678
$ws.Range("H124").Value = 8976.23
$ws.Range("I124").Value = 7353.6
$ws.Range("J124").Value = 9990.375
$ws.Range("K124").Value = 22060.8
$ws.Range("L124").Value = 29971.125
$ws.Range("M124").Value = -17150.8
$ws.Range("N124").Value = -39791.125
$ws.Range("H131").Value = 56981.61
$ws.Range("I131").Value = 84127.414
$ws.Range("J131").Value = 2690
$ws.Range("K131").Value = 252382.242
$ws.Range("L131").Value = 8070
$ws.Range("M131").Value = -247342.242
$ws.Range("N131").Value = -18150
$ws.Range("H137").Value = 2764.3845
$ws.Range("I137").Value = 1870
$ws.Range("J137").Value = 4195.4
$ws.Range("K137").Value = 5610
$ws.Range("L137").Value = 12586.2
$ws.Range("M137").Value = -510
$ws.Range("N137").Value = -22786.2

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 3118.2
$ws.Range("I132").Value = 2250.3333
$ws.Range("K132").Value = 6750.999899999999
$ws.Range("M132").Value = -4220.999899999999

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 354.57895
$ws.Range("I55").Value = 356.36365
$ws.Range("J55").Value = 352.125
$ws.Range("K55").Value = 356.36365
$ws.Range("L55").Value = 352.125
$ws.Range("M55").Value = -183.36365
$ws.Range("N55").Value = -698.125
$ws.Range("H93").Value = 2805.6428
$ws.Range("I93").Value = 2805.6428
$ws.Range("J93").Value = 0
$ws.Range("K93").Value = 2805.6428
$ws.Range("L93").Value = 0
$ws.Range("M93").ClearContents()
$ws.Range("N93").Value = -1557.6428
$ws.Range("H132").Value = 4582.6113
$ws.Range("I132").Value = 3598.7144
$ws.Range("K132").Value = 10796.1432
$ws.Range("M132").Value = -8266.143199999999

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 3949
$ws.Range("I132").Value = 3685.963
$ws.Range("K132").Value = 11057.889
$ws.Range("M132").Value = -8527.889000000001
